$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 30000
$ws.Range("J68").Value = 30000
$ws.Range("L68").Value = 30000
$ws.Range("N68").Value = -31498
$ws.Range("H71").Value = 30000
$ws.Range("J71").Value = 30000
$ws.Range("L71").Value = 90000
$ws.Range("N71").Value = -97488
$ws.Range("H74").Value = 4246.7334
$ws.Range("I74").Value = 3744.3333
$ws.Range("K74").Value = 3744.3333
$ws.Range("M74").Value = -2808.3333
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H76").Value = 3196.4211
$ws.Range("I76").Value = 3133.2222
$ws.Range("J76").Value = 3253.3
$ws.Range("K76").Value = 3133.2222
$ws.Range("L76").Value = 3253.3
$ws.Range("M76").Value = -2818.2222
$ws.Range("N76").Value = -3883.3
$ws.Range("H77").Value = 4246.7334
$ws.Range("I77").Value = 3744.3333
$ws.Range("K77").Value = 18721.6665
$ws.Range("M77").Value = -14041.6665
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H79").Value = 3196.4211
$ws.Range("I79").Value = 3133.2222
$ws.Range("J79").Value = 3253.3
$ws.Range("K79").Value = 3133.2222
$ws.Range("L79").Value = 3253.3
$ws.Range("M79").Value = -2041.2222
$ws.Range("N79").Value = -5437.3
$ws.Range("H93").Value = 36111.06
$ws.Range("J93").Value = 36111.06
$ws.Range("L93").Value = 36111.06
$ws.Range("N93").Value = -41103.06
$ws.Range("H129").Value = 1593.098
$ws.Range("I129").Value = 363.45456
$ws.Range("J129").Value = 1931.25
$ws.Range("K129").Value = 1090.36368
$ws.Range("L129").Value = 5793.75
$ws.Range("M129").Value = 3909.63632
$ws.Range("N129").Value = -15793.75
$ws.Range("H137").Value = 28572530
$ws.Range("I137").Value = 941.2
$ws.Range("K137").Value = 2823.6
$ws.Range("M137").Value = -273.6000000000004
$ws.Range("H141").Value = 3026.0625
$ws.Range("I141").Value = 1126.8125
$ws.Range("J141").Value = 4925.3125
$ws.Range("K141").Value = 3380.4375
$ws.Range("L141").Value = 14775.9375
$ws.Range("M141").Value = 1799.5625
$ws.Range("N141").Value = -25135.9375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 1700
$ws.Range("I12").Value = 1700
$ws.Range("K12").Value = 1700
$ws.Range("M12").Value = -1527
$ws.Range("H45").Value = 1265.9166
$ws.Range("I45").Value = 1185.4667
$ws.Range("J45").Value = 1400
$ws.Range("K45").Value = 1185.4667
$ws.Range("L45").Value = 1400
$ws.Range("M45").Value = -808.4666999999999
$ws.Range("N45").Value = -2154
$ws.Range("H74").Value = 33336736
$ws.Range("I74").Value = 62501910
$ws.Range("K74").Value = 62501910
$ws.Range("M74").Value = -62501036
$ws.Range("H77").Value = 33336736
$ws.Range("I77").Value = 62501910
$ws.Range("K77").Value = 312509550
$ws.Range("M77").Value = -312505182
$ws.Range("H110").Value = 640
$ws.Range("I110").Value = 555
$ws.Range("J110").Value = 725
$ws.Range("K110").Value = 555
$ws.Range("L110").Value = 725
$ws.Range("M110").Value = 1490
$ws.Range("N110").Value = -4815
$ws.Range("H117").Value = 24029.54
$ws.Range("J117").Value = 24029.54
$ws.Range("L117").Value = 24029.54
$ws.Range("N117").Value = -33207.54

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2997.0417
$ws.Range("I105").Value = 1533.9131
$ws.Range("J105").Value = 4343.12
$ws.Range("K105").Value = 1533.9131
$ws.Range("L105").Value = 4343.12
$ws.Range("M105").Value = 213.0869
$ws.Range("N105").Value = -7837.12
$ws.Range("H107").Value = 861.9375
$ws.Range("I107").Value = 863.6429000000001
$ws.Range("J107").Value = 850
$ws.Range("K107").Value = 863.6429000000001
$ws.Range("L107").Value = 850
$ws.Range("M107").Value = 1056.3571
$ws.Range("N107").Value = -4690
$ws.Range("H118").Value = 7692.5
$ws.Range("J118").Value = 7692.5
$ws.Range("L118").Value = 7692.5
$ws.Range("N118").Value = -11006.5
$ws.Range("H134").Value = 5662.222
$ws.Range("I134").Value = 4173
$ws.Range("J134").Value = 6087.7144
$ws.Range("K134").Value = 12519
$ws.Range("L134").Value = 18263.1432
$ws.Range("M134").Value = -9984
$ws.Range("N134").Value = -23333.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 20319.5
$ws.Range("J51").Value = 20319.5
$ws.Range("L51").Value = 20319.5
$ws.Range("N51").Value = -21337.5
$ws.Range("H80").Value = 17546670
$ws.Range("I80").Value = 37039416
$ws.Range("J80").Value = 3197.9
$ws.Range("K80").Value = 37039416
$ws.Range("L80").Value = 3197.9
$ws.Range("M80").Value = -37038418
$ws.Range("N80").Value = -5193.9
$ws.Range("H83").Value = 17546670
$ws.Range("I83").Value = 37039416
$ws.Range("J83").Value = 3197.9
$ws.Range("K83").Value = 185197080
$ws.Range("L83").Value = 15989.5
$ws.Range("M83").Value = -185192088
$ws.Range("N83").Value = -25973.5
$ws.Range("H102").Value = 7555
$ws.Range("I102").Value = 9428.888999999999
$ws.Range("J102").Value = 1933.3334
$ws.Range("K102").Value = 9428.888999999999
$ws.Range("L102").Value = 1933.3334
$ws.Range("M102").Value = -7806.888999999999
$ws.Range("N102").Value = -5177.3334
$ws.Range("H126").Value = 5274.1113
$ws.Range("I126").Value = 2891.75
$ws.Range("J126").Value = 7180
$ws.Range("K126").Value = 8675.25
$ws.Range("L126").Value = 21540
$ws.Range("M126").Value = -6205.25
$ws.Range("N126").Value = -26480
$ws.Range("H132").Value = 6600.9062
$ws.Range("I132").Value = 6732.069
$ws.Range("K132").Value = 20196.207
$ws.Range("M132").Value = -17666.207

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4949.28
$ws.Range("I40").Value = 7466.375
$ws.Range("J40").Value = 3764.7646
$ws.Range("K40").Value = 7466.375
$ws.Range("L40").Value = 3764.7646
$ws.Range("M40").Value = -7330.375
$ws.Range("N40").Value = -4036.7646
$ws.Range("H122").Value = 6580.44
$ws.Range("I122").Value = 7680.9287
$ws.Range("J122").Value = 5179.8184
$ws.Range("K122").Value = 23042.7861
$ws.Range("L122").Value = 15539.4552
$ws.Range("M122").Value = -20592.7861
$ws.Range("N122").Value = -20439.4552

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3450.5
$ws.Range("I81").Value = 2001
$ws.Range("J81").Value = 4900
$ws.Range("K81").Value = 4002
$ws.Range("L81").Value = 9800
$ws.Range("M81").Value = -2941
$ws.Range("N81").Value = -11922
$ws.Range("H84").Value = 3450.5
$ws.Range("I84").Value = 2001
$ws.Range("J84").Value = 4900
$ws.Range("K84").Value = 20010
$ws.Range("L84").Value = 49000
$ws.Range("M84").Value = -14706
$ws.Range("N84").Value = -59608

